$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30000
$ws.Range("I3").Value = 30000
$ws.Range("K3").Value = 30000
$ws.Range("M3").Value = -29886
$ws.Range("H17").Value = 2706
$ws.Range("I17").Value = 1440
$ws.Range("J17").Value = 2959.2
$ws.Range("K17").Value = 4320
$ws.Range("L17").Value = 8877.599999999999
$ws.Range("M17").Value = -4152
$ws.Range("N17").Value = -9213.599999999999
$ws.Range("H64").Value = 333336480
$ws.Range("I64").Value = 4724.5
$ws.Range("J64").Value = 1000000000
$ws.Range("K64").Value = 4724.5
$ws.Range("L64").Value = 1000000000
$ws.Range("M64").Value = -4476.5
$ws.Range("N64").Value = -1000000496
$ws.Range("H67").Value = 333336480
$ws.Range("I67").Value = 4724.5
$ws.Range("J67").Value = 1000000000
$ws.Range("K67").Value = 4724.5
$ws.Range("L67").Value = 1000000000
$ws.Range("M67").Value = -3866.5
$ws.Range("N67").Value = -1000001716
$ws.Range("H100").Value = 1229.5714
$ws.Range("I100").Value = 720.6
$ws.Range("J100").Value = 2502
$ws.Range("K100").Value = 720.6
$ws.Range("L100").Value = 2502
$ws.Range("M100").Value = -179.6
$ws.Range("N100").Value = -3584
$ws.Range("H102").Value = 30000
$ws.Range("I102").Value = 30000
$ws.Range("K102").Value = 30000
$ws.Range("M102").Value = -26755
$ws.Range("H106").Value = 3481.8
$ws.Range("I106").Value = 3127.25
$ws.Range("K106").Value = 3127.25
$ws.Range("M106").Value = -2496.25
$ws.Range("H131").Value = 835516.3
$ws.Range("I131").Value = 1112525
$ws.Range("K131").Value = 3337575
$ws.Range("M131").Value = -3332535
$ws.Range("H138").Value = 3300.64
$ws.Range("I138").Value = 2736.2942
$ws.Range("K138").Value = 8208.882599999999
$ws.Range("M138").Value = -3068.882599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3531.625
$ws.Range("I32").Value = 3952.7144
$ws.Range("J32").Value = 584
$ws.Range("K32").Value = 3952.7144
$ws.Range("L32").Value = 584
$ws.Range("M32").Value = -3665.7144
$ws.Range("N32").Value = -1158
$ws.Range("H61").Value = 2505.05
$ws.Range("I61").Value = 1700.1666
$ws.Range("K61").Value = 1700.1666
$ws.Range("M61").Value = -1488.1666
$ws.Range("H74").Value = 114755.18
$ws.Range("I74").Value = 135667.97
$ws.Range("K74").Value = 135667.97
$ws.Range("M74").Value = -134793.97
$ws.Range("H77").Value = 114755.18
$ws.Range("I77").Value = 135667.97
$ws.Range("K77").Value = 678339.85
$ws.Range("M77").Value = -673971.85
$ws.Range("H122").Value = 1864.4166
$ws.Range("I122").Value = 1486.2222
$ws.Range("K122").Value = 4458.6666
$ws.Range("M122").Value = -2008.6666
$ws.Range("H132").Value = 3222.5293
$ws.Range("I132").Value = 2948.5833
$ws.Range("J132").Value = 3880
$ws.Range("K132").Value = 8845.749899999999
$ws.Range("L132").Value = 11640
$ws.Range("M132").Value = -6315.749899999999
$ws.Range("N132").Value = -16700
$ws.Range("H136").Value = 2505.05
$ws.Range("I136").Value = 1700.1666
$ws.Range("K136").Value = 5100.4998
$ws.Range("M136").Value = -2550.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 62500464
$ws.Range("J80").Value = 510.25
$ws.Range("L80").Value = 510.25
$ws.Range("N80").Value = -2506.25
$ws.Range("H83").Value = 62500464
$ws.Range("J83").Value = 510.25
$ws.Range("L83").Value = 2551.25
$ws.Range("N83").Value = -12535.25
$ws.Range("H86").Value = 2171.7917
$ws.Range("I86").Value = 1554.1333
$ws.Range("K86").Value = 1554.1333
$ws.Range("M86").Value = -431.1333
$ws.Range("H89").Value = 2171.7917
$ws.Range("I89").Value = 1554.1333
$ws.Range("K89").Value = 7770.666499999999
$ws.Range("M89").Value = -2154.666499999999
$ws.Range("H94").Value = 86957840
$ws.Range("I94").Value = 95239200
$ws.Range("K94").Value = 95239200
$ws.Range("M94").Value = -95238749
$ws.Range("H134").Value = 3204.5833
$ws.Range("I134").Value = 2556
$ws.Range("J134").Value = 4501.75
$ws.Range("K134").Value = 7668
$ws.Range("L134").Value = 13505.25
$ws.Range("M134").Value = -5133
$ws.Range("N134").Value = -18575.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2774.2974
$ws.Range("I31").Value = 1786.9
$ws.Range("J31").Value = 7006
$ws.Range("K31").Value = 1786.9
$ws.Range("L31").Value = 7006
$ws.Range("M31").Value = -1491.9
$ws.Range("N31").Value = -7596
$ws.Range("H34").Value = 2774.2974
$ws.Range("I34").Value = 1786.9
$ws.Range("J34").Value = 7006
$ws.Range("K34").Value = 1786.9
$ws.Range("L34").Value = 7006
$ws.Range("M34").Value = -1584.9
$ws.Range("N34").Value = -7410
$ws.Range("H58").Value = 2170.6667
$ws.Range("I58").Value = 1256
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 1256
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -1053
$ws.Range("N58").Value = -4406
$ws.Range("H69").Value = 11111
$ws.Range("I69").Value = 11111
$ws.Range("K69").Value = 11111
$ws.Range("M69").Value = -10362
$ws.Range("H72").Value = 11111
$ws.Range("I72").Value = 11111
$ws.Range("K72").Value = 33333
$ws.Range("M72").Value = -29589
$ws.Range("H132").Value = 11909947
$ws.Range("I132").Value = 4231
$ws.Range("K132").Value = 12693
$ws.Range("M132").Value = -10163
$ws.Range("H134").Value = 2225.275
$ws.Range("I134").Value = 2124.8125
$ws.Range("J134").Value = 2627.125
$ws.Range("K134").Value = 6374.4375
$ws.Range("L134").Value = 7881.375
$ws.Range("M134").Value = -3839.4375
$ws.Range("N134").Value = -12951.375
$ws.Range("H136").Value = 2170.6667
$ws.Range("I136").Value = 1256
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3768
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1218
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1596.4736
$ws.Range("J122").Value = 1674.0555
$ws.Range("L122").Value = 15066.4995
$ws.Range("N122").Value = -19966.4995
$ws.Range("H132").Value = 3557.1428
$ws.Range("I132").Value = 900
$ws.Range("K132").Value = 8100
$ws.Range("M132").Value = -5570
$ws.Range("H139").Value = 3411.3076
$ws.Range("I139").Value = 2454
$ws.Range("K139").Value = 7362
$ws.Range("M139").Value = -2222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3268.2144
$ws.Range("I122").Value = 2341.8572
$ws.Range("K122").Value = 7025.571599999999
$ws.Range("M122").Value = -4575.571599999999
$ws.Range("H132").Value = 3686.3572
$ws.Range("I132").Value = 2701.125
$ws.Range("K132").Value = 8103.375
$ws.Range("M132").Value = -5573.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 33347504
$ws.Range("I2").Value = 5002.5
$ws.Range("J2").Value = 50018750
$ws.Range("K2").Value = 5002.5
$ws.Range("L2").Value = 50018750
$ws.Range("M2").Value = -4890.5
$ws.Range("N2").Value = -50018974
$ws.Range("H40").Value = 61792.707
$ws.Range("I40").Value = 61792.707
$ws.Range("K40").Value = 61792.707
$ws.Range("M40").Value = -61656.707
$ws.Range("H93").Value = 1473
$ws.Range("I93").Value = 1532.1538
$ws.Range("K93").Value = 1532.1538
$ws.Range("M93").Value = -284.1538
$ws.Range("H100").Value = 6583
$ws.Range("I100").Value = 6500
$ws.Range("K100").Value = 6500
$ws.Range("M100").Value = -5959

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6319.857
$ws.Range("I62").Value = 6068.6
$ws.Range("K62").Value = 6068.6
$ws.Range("M62").Value = -5444.6
$ws.Range("H65").Value = 6319.857
$ws.Range("I65").Value = 6068.6
$ws.Range("K65").Value = 30343
$ws.Range("M65").Value = -27223
$ws.Range("H101").Value = 44409.25
$ws.Range("J101").Value = 44409.25
$ws.Range("L101").Value = 44409.25
$ws.Range("N101").Value = -50899.25
$ws.Range("H107").Value = 579
$ws.Range("I107").Value = 316.9091
$ws.Range("K107").Value = 950.7273
$ws.Range("M107").Value = 969.2727
$ws.Range("H126").Value = 2559.3333
$ws.Range("I126").Value = 2476.5
$ws.Range("J126").Value = 2725
$ws.Range("K126").Value = 7429.5
$ws.Range("L126").Value = 8175
$ws.Range("M126").Value = -4959.5
$ws.Range("N126").Value = -13115
$ws.Range("H132").Value = 4286.4346
$ws.Range("I132").Value = 3647.1428
$ws.Range("K132").Value = 10941.4284
$ws.Range("M132").Value = -8411.428400000001
$ws.Range("H135").Value = 123388.44
$ws.Range("J135").Value = 123388.44
$ws.Range("L135").Value = 123388.44
$ws.Range("N135").Value = -133528.44
